$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Header 2 is the document's primary header (the one that actually holds
# content here) -- it contains the BTec logo inline picture, currently
# named "image1.jpg". Rename it to "image2.jpg".
$header = $sec.Headers.Item(2)
$btecLogo = $header.Range.InlineShapes.Item(1)
$btecLogo.Name = "image2.jpg"

# Footer 1 holds a Pearson logo inline picture currently named
# "image2.png". Rename it to "image1.png".
$footer1 = $sec.Footers.Item(1)
$pearsonLogo1 = $footer1.Range.InlineShapes.Item(1)
$pearsonLogo1.Name = "image1.png"

# Footer 2 holds the matching Pearson logo inline picture, also currently
# named "image2.png". Rename it to "image1.png" as well.
$footer2 = $sec.Footers.Item(2)
$pearsonLogo2 = $footer2.Range.InlineShapes.Item(1)
$pearsonLogo2.Name = "image1.png"

Write-Host "Renamed BTec logo to image2.jpg; Pearson logos (footer1 & footer2) to image1.png."
